$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-20 and clear/replace rows 21-29, then add new rows 30-43 ---
# row 2
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "0x63d4e7e"
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = "0x9a"
$ws.Cells.Item(2,5).Value = 23

# row 3
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "0x55"
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = "0x9f"
$ws.Cells.Item(3,5).Value = 102

# row 4
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "0xe1"
$ws.Cells.Item(4,3).Value = 120
$ws.Cells.Item(4,4).Value = "0x48"
$ws.Cells.Item(4,5).Value = 434

# row 5
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "0xf1"
$ws.Cells.Item(5,3).Value = 425
$ws.Cells.Item(5,4).Value = "0x108"
$ws.Cells.Item(5,5).Value = 1

# row 6
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "0x101"
$ws.Cells.Item(6,3).Value = 10
$ws.Cells.Item(6,4).Value = "0x3c"
$ws.Cells.Item(6,5).Value = 9

# row 7
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "0x8c"
$ws.Cells.Item(7,3).Value = 10
$ws.Cells.Item(7,4).Value = "0xc7"
$ws.Cells.Item(7,5).Value = 44

# row 8
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "0x111"
$ws.Cells.Item(8,3).Value = 24
$ws.Cells.Item(8,4).Value = "0xd0"
$ws.Cells.Item(8,5).Value = 44

# row 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "0x151"
$ws.Cells.Item(9,3).Value = 8
$ws.Cells.Item(9,4).Value = "0xd1"
$ws.Cells.Item(9,5).Value = 44

# row 10
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "0x1c"
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = "0xd2"
$ws.Cells.Item(10,5).Value = 44

# row 11
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "0x121"
$ws.Cells.Item(11,3).Value = 6
$ws.Cells.Item(11,4).Value = "0xd3"
$ws.Cells.Item(11,5).Value = 44

# row 12
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "0xa8"
$ws.Cells.Item(12,3).Value = 9
$ws.Cells.Item(12,4).Value = "0xd4"
$ws.Cells.Item(12,5).Value = 44

# row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "0xb0"
$ws.Cells.Item(13,3).Value = 7
$ws.Cells.Item(13,4).Value = "0x170"
$ws.Cells.Item(13,5).Value = 1

# row 14
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "0xb1"
$ws.Cells.Item(14,3).Value = 8
$ws.Cells.Item(14,4).Value = "0x202"
$ws.Cells.Item(14,5).Value = 9

# row 15
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "0xcc"
$ws.Cells.Item(15,3).Value = 8
$ws.Cells.Item(15,4).Value = "0x203"
$ws.Cells.Item(15,5).Value = 9

# row 16
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "0xcd"
$ws.Cells.Item(16,3).Value = 6
$ws.Cells.Item(16,4).Value = "0x205"
$ws.Cells.Item(16,5).Value = 44

# row 17
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "0x128"
$ws.Cells.Item(17,3).Value = 5
$ws.Cells.Item(17,4).Value = "0x220"
$ws.Cells.Item(17,5).Value = 44

# row 18
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "0x140"
$ws.Cells.Item(18,3).Value = 5
$ws.Cells.Item(18,4).Value = "0x221"
$ws.Cells.Item(18,5).Value = 44

# row 19
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "0x1d0"
$ws.Cells.Item(19,3).Value = 3
$ws.Cells.Item(19,4).Value = "0x72"
$ws.Cells.Item(19,5).Value = 1

# row 20
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "0x2aa"
$ws.Cells.Item(20,3).Value = 3
$ws.Cells.Item(20,4).Value = "0x73"
$ws.Cells.Item(20,5).Value = 1

# row 21
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).ClearContents()
$ws.Cells.Item(21,3).ClearContents()
$ws.Cells.Item(21,4).Value = "0xa0"
$ws.Cells.Item(21,5).Value = 3

# row 22
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).ClearContents()
$ws.Cells.Item(22,3).ClearContents()
$ws.Cells.Item(22,4).Value = "0xa1"
$ws.Cells.Item(22,5).Value = 3

# row 23
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).ClearContents()
$ws.Cells.Item(23,3).ClearContents()
$ws.Cells.Item(23,4).Value = "0xa2"
$ws.Cells.Item(23,5).Value = 3

# row 24
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).ClearContents()
$ws.Cells.Item(24,3).ClearContents()
$ws.Cells.Item(24,4).Value = "0xa3"
$ws.Cells.Item(24,5).Value = 3

# row 25
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).ClearContents()
$ws.Cells.Item(25,3).ClearContents()
$ws.Cells.Item(25,4).Value = "0xa4"
$ws.Cells.Item(25,5).Value = 3

# row 26
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).ClearContents()
$ws.Cells.Item(26,3).ClearContents()
$ws.Cells.Item(26,4).Value = "0xa6"
$ws.Cells.Item(26,5).Value = 3

# row 27
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).ClearContents()
$ws.Cells.Item(27,3).ClearContents()
$ws.Cells.Item(27,4).Value = "0xbc"
$ws.Cells.Item(27,5).Value = 3

# row 28
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).ClearContents()
$ws.Cells.Item(28,3).ClearContents()
$ws.Cells.Item(28,4).Value = "0x2c0"
$ws.Cells.Item(28,5).Value = 3

# row 29
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).ClearContents()
$ws.Cells.Item(29,3).ClearContents()
$ws.Cells.Item(29,4).Value = "0x2c1"
$ws.Cells.Item(29,5).Value = 3

# row 30
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).ClearContents()
$ws.Cells.Item(30,3).ClearContents()
$ws.Cells.Item(30,4).Value = "0x59"
$ws.Cells.Item(30,5).Value = 34

# row 31
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).ClearContents()
$ws.Cells.Item(31,3).ClearContents()
$ws.Cells.Item(31,4).Value = "0x131"
$ws.Cells.Item(31,5).Value = 1

# row 32
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).ClearContents()
$ws.Cells.Item(32,3).ClearContents()
$ws.Cells.Item(32,4).Value = "0xbd"
$ws.Cells.Item(32,5).Value = 4

# row 33
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).ClearContents()
$ws.Cells.Item(33,3).ClearContents()
$ws.Cells.Item(33,4).Value = "0xbe"
$ws.Cells.Item(33,5).Value = 4

# row 34
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).ClearContents()
$ws.Cells.Item(34,3).ClearContents()
$ws.Cells.Item(34,4).Value = "0xbf"
$ws.Cells.Item(34,5).Value = 4

# row 35
$ws.Cells.Item(35,1).Value = 33
$ws.Cells.Item(35,2).ClearContents()
$ws.Cells.Item(35,3).ClearContents()
$ws.Cells.Item(35,4).Value = "0x171"
$ws.Cells.Item(35,5).Value = 1

# row 36
$ws.Cells.Item(36,1).Value = 34
$ws.Cells.Item(36,2).ClearContents()
$ws.Cells.Item(36,3).ClearContents()
$ws.Cells.Item(36,4).Value = "0xc0"
$ws.Cells.Item(36,5).Value = 11

# row 37
$ws.Cells.Item(37,1).Value = 35
$ws.Cells.Item(37,2).ClearContents()
$ws.Cells.Item(37,3).ClearContents()
$ws.Cells.Item(37,4).Value = "0x9b"
$ws.Cells.Item(37,5).Value = 18

# row 38
$ws.Cells.Item(38,1).Value = 36
$ws.Cells.Item(38,2).ClearContents()
$ws.Cells.Item(38,3).ClearContents()
$ws.Cells.Item(38,4).Value = "0x211"
$ws.Cells.Item(38,5).Value = 7

# row 39
$ws.Cells.Item(39,1).Value = 37
$ws.Cells.Item(39,2).ClearContents()
$ws.Cells.Item(39,3).ClearContents()
$ws.Cells.Item(39,4).Value = "0x210"
$ws.Cells.Item(39,5).Value = 4

# row 40
$ws.Cells.Item(40,1).Value = 38
$ws.Cells.Item(40,2).ClearContents()
$ws.Cells.Item(40,3).ClearContents()
$ws.Cells.Item(40,4).Value = "0x37"
$ws.Cells.Item(40,5).Value = 33

# row 41
$ws.Cells.Item(41,1).Value = 39
$ws.Cells.Item(41,2).ClearContents()
$ws.Cells.Item(41,3).ClearContents()
$ws.Cells.Item(41,4).Value = "0x3b"
$ws.Cells.Item(41,5).Value = 7

# row 42
$ws.Cells.Item(42,1).Value = 40
$ws.Cells.Item(42,2).ClearContents()
$ws.Cells.Item(42,3).ClearContents()
$ws.Cells.Item(42,4).Value = "0x4a"
$ws.Cells.Item(42,5).Value = 200

# row 43
$ws.Cells.Item(43,1).Value = 41
$ws.Cells.Item(43,2).ClearContents()
$ws.Cells.Item(43,3).ClearContents()
$ws.Cells.Item(43,4).Value = "0x30"
$ws.Cells.Item(43,5).Value = 1

# --- Apply matching style (bold/centered/bordered) to new A-column cells (rows 30-43) ---
$ws.Range("A29").Copy()
$ws.Range("A30:A43").PasteSpecial(-4122)
$excel.CutCopyMode = 0
